# Applies the "Deploying to gh-pages" content refresh to the FHIR
# StructureDefinition spreadsheet:
#   - Metadata sheet: bump Version/Date, replace the duplicated "Contact"
#     row with Publisher/Jurisdiction details, and drop the extra duplicate
#     row (net: 21 -> 20 rows).
#   - Elements sheet: give the root Extension row its real Short/Definition
#     text instead of the generic "Extension" / "An Extension" placeholders.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Remove the second of the two duplicated "Contact" / "No display for
# ContactDetail" rows (originally row 11); everything below shifts up.
$meta.Rows(11).Delete()

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Updated publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes a Jurisdiction row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Give the root Extension element its real Short/Definition text
$elements.Range("K2").Value = "Vision Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member has vision benefit coverage: Y or N"
